$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 17 new rows before row 926 to make room for the full, alphabetized
# region list (ISS update 2021-04-02). Existing rows 926-929 (Campania,
# Lombardia, Piemonte, Toscana) shift down to 930, 935, 938, 942.
$ws.Rows("926:942").Insert()

$ws.Cells.Item(926, 1).Value = 'Abruzzo'
$ws.Cells.Item(926, 2).Value = 0.83
$ws.Cells.Item(926, 3).Value = 0.72
$ws.Cells.Item(926, 4).Value = 0.93
$ws.Cells.Item(926, 5).Value = 44265
$ws.Cells.Item(926, 6).Value = 44278
$ws.Cells.Item(927, 1).Value = 'Basilicata'
$ws.Cells.Item(927, 2).Value = 1.07
$ws.Cells.Item(927, 3).Value = 0.6899999999999999
$ws.Cells.Item(927, 4).Value = 1.52
$ws.Cells.Item(927, 5).Value = 44265
$ws.Cells.Item(927, 6).Value = 44278
$ws.Cells.Item(928, 1).Value = 'Bolzano'
$ws.Cells.Item(928, 2).Value = 0.8100000000000001
$ws.Cells.Item(928, 3).Value = 0.6899999999999999
$ws.Cells.Item(928, 4).Value = 1.02
$ws.Cells.Item(928, 5).Value = 44265
$ws.Cells.Item(928, 6).Value = 44278
$ws.Cells.Item(929, 1).Value = 'Calabria'
$ws.Cells.Item(929, 2).Value = 1.29
$ws.Cells.Item(929, 3).Value = 0.85
$ws.Cells.Item(929, 4).Value = 1.71
$ws.Cells.Item(929, 5).Value = 44265
$ws.Cells.Item(929, 6).Value = 44278
$ws.Cells.Item(930, 1).Value = 'Campania'
$ws.Cells.Item(930, 2).Value = 1.22
$ws.Cells.Item(930, 3).Value = 0.96
$ws.Cells.Item(930, 4).Value = 1.46
$ws.Cells.Item(930, 5).Value = 44265
$ws.Cells.Item(930, 6).Value = 44278
$ws.Cells.Item(931, 1).Value = 'Emilia R.'
$ws.Cells.Item(931, 2).Value = 0.87
$ws.Cells.Item(931, 3).Value = 0.78
$ws.Cells.Item(931, 4).Value = 1.02
$ws.Cells.Item(931, 5).Value = 44265
$ws.Cells.Item(931, 6).Value = 44278
$ws.Cells.Item(932, 1).Value = 'Friuli V.G.'
$ws.Cells.Item(932, 2).Value = 1.01
$ws.Cells.Item(932, 3).Value = 0.72
$ws.Cells.Item(932, 4).Value = 1.31
$ws.Cells.Item(932, 5).Value = 44265
$ws.Cells.Item(932, 6).Value = 44278
$ws.Cells.Item(933, 1).Value = 'Lazio'
$ws.Cells.Item(933, 2).Value = 0.96
$ws.Cells.Item(933, 3).Value = 0.83
$ws.Cells.Item(933, 4).Value = 1.04
$ws.Cells.Item(933, 5).Value = 44265
$ws.Cells.Item(933, 6).Value = 44278
$ws.Cells.Item(934, 1).Value = 'Liguria'
$ws.Cells.Item(934, 2).Value = 1.03
$ws.Cells.Item(934, 3).Value = 0.95
$ws.Cells.Item(934, 4).Value = 1.14
$ws.Cells.Item(934, 5).Value = 44265
$ws.Cells.Item(934, 6).Value = 44278
$ws.Cells.Item(935, 1).Value = 'Lombardia'
$ws.Cells.Item(935, 2).Value = 0.92
$ws.Cells.Item(935, 3).Value = 0.72
$ws.Cells.Item(935, 4).Value = 1.13
$ws.Cells.Item(935, 5).Value = 44265
$ws.Cells.Item(935, 6).Value = 44278
$ws.Cells.Item(936, 1).Value = 'Marche'
$ws.Cells.Item(936, 2).Value = 1.03
$ws.Cells.Item(936, 3).Value = 0.77
$ws.Cells.Item(936, 4).Value = 1.26
$ws.Cells.Item(936, 5).Value = 44265
$ws.Cells.Item(936, 6).Value = 44278
$ws.Cells.Item(937, 1).Value = 'Molise'
$ws.Cells.Item(937, 2).Value = 1.04
$ws.Cells.Item(937, 3).Value = 0.5
$ws.Cells.Item(937, 4).Value = 1.58
$ws.Cells.Item(937, 5).Value = 44265
$ws.Cells.Item(937, 6).Value = 44278
$ws.Cells.Item(938, 1).Value = 'Piemonte'
$ws.Cells.Item(938, 2).Value = 0.98
$ws.Cells.Item(938, 3).Value = 0.74
$ws.Cells.Item(938, 4).Value = 1.22
$ws.Cells.Item(938, 5).Value = 44265
$ws.Cells.Item(938, 6).Value = 44278
$ws.Cells.Item(939, 1).Value = 'Puglia'
$ws.Cells.Item(939, 2).Value = 1.1
$ws.Cells.Item(939, 3).Value = 0.95
$ws.Cells.Item(939, 4).Value = 1.23
$ws.Cells.Item(939, 5).Value = 44265
$ws.Cells.Item(939, 6).Value = 44278
$ws.Cells.Item(940, 1).Value = 'Sardegna'
$ws.Cells.Item(940, 2).Value = 1.18
$ws.Cells.Item(940, 3).Value = 0.92
$ws.Cells.Item(940, 4).Value = 1.54
$ws.Cells.Item(940, 5).Value = 44265
$ws.Cells.Item(940, 6).Value = 44278
$ws.Cells.Item(941, 1).Value = 'Sicilia'
$ws.Cells.Item(941, 2).Value = 1.09
$ws.Cells.Item(941, 3).Value = 0.98
$ws.Cells.Item(941, 4).Value = 1.21
$ws.Cells.Item(941, 5).Value = 44265
$ws.Cells.Item(941, 6).Value = 44278
$ws.Cells.Item(942, 1).Value = 'Toscana'
$ws.Cells.Item(942, 2).Value = 1.06
$ws.Cells.Item(942, 3).Value = 1
$ws.Cells.Item(942, 4).Value = 1.11
$ws.Cells.Item(942, 5).Value = 44265
$ws.Cells.Item(942, 6).Value = 44278
$ws.Cells.Item(943, 1).Value = 'Trento'
$ws.Cells.Item(943, 2).Value = 0.89
$ws.Cells.Item(943, 3).Value = 0.78
$ws.Cells.Item(943, 4).Value = 1.08
$ws.Cells.Item(943, 5).Value = 44265
$ws.Cells.Item(943, 6).Value = 44278
$ws.Cells.Item(944, 1).Value = 'Umbria'
$ws.Cells.Item(944, 2).Value = 0.88
$ws.Cells.Item(944, 3).Value = 0.75
$ws.Cells.Item(944, 4).Value = 1.03
$ws.Cells.Item(944, 5).Value = 44265
$ws.Cells.Item(944, 6).Value = 44278
$ws.Cells.Item(945, 1).Value = 'Valle d''Aosta'
$ws.Cells.Item(945, 2).Value = 1.61
$ws.Cells.Item(945, 3).Value = 1.22
$ws.Cells.Item(945, 4).Value = 2.2
$ws.Cells.Item(945, 5).Value = 44265
$ws.Cells.Item(945, 6).Value = 44278
$ws.Cells.Item(946, 1).Value = 'Veneto'
$ws.Cells.Item(946, 2).Value = 1.11
$ws.Cells.Item(946, 3).Value = 0.9
$ws.Cells.Item(946, 4).Value = 1.28
$ws.Cells.Item(946, 5).Value = 44265
$ws.Cells.Item(946, 6).Value = 44278
